$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the first
#    paragraph (the Heading1 title "Play Buffalo Blitz Slot for Free -
#    Review and Gameplay"). The new paragraph has:
#      - an empty leading run (<w:r/>) matching the document's existing
#        paragraph pattern
#      - a bold run containing "Meta description"
#      - a normal run containing the rest of the description text
# -----------------------------------------------------------------------

$metaText = ": Get ready to play Buffalo Blitz slot for free! Read our review to learn more about its features, gameplay mechanics, and unique paylines."

# First create a brand-new, empty paragraph directly after paragraph 1
# (InsertXML replaces the contents of whatever paragraph it lands on, so
# we must target a fresh, empty paragraph rather than collapsing a range
# that still belongs to the following "GAME MECHANICS AND PAYLINES"
# paragraph - otherwise that paragraph's text would be overwritten).
$firstPara = $d.Paragraphs.Item(1)
$endOfFirst = $firstPara.Range
$endOfFirst.Collapse(0)  # wdCollapseEnd
$endOfFirst.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$insertPoint = $newPara.Range
$insertPoint.Collapse(1)  # wdCollapseStart

$metaXml = "<?xml version='1.0' encoding='UTF-8' standalone='yes'?>" +
  "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
  "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
  "<pkg:xmlData>" +
  "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
  "<w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>" + $metaText + "</w:t></w:r></w:p></w:body>" +
  "</w:document></pkg:xmlData></pkg:part></pkg:package>"

[void]$insertPoint.InsertXML($metaXml)

# -----------------------------------------------------------------------
# 2. Remove the duplicate bold "Play Buffalo Blitz Slot for Free - Review
#    and Gameplay" paragraph that used to sit near the end of the document
#    (right before the italic meta-description-like paragraph).
# -----------------------------------------------------------------------

$count = $d.Paragraphs.Count
$boldPara = $d.Paragraphs.Item($count - 1)
$boldPara.Range.Delete()

# -----------------------------------------------------------------------
# 3. Replace the text of the final (italic) paragraph with the new image
#    generation prompt, keeping its italic run formatting intact.
# -----------------------------------------------------------------------

$newCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($newCount)
$start = $lastPara.Range.Start
$end = $lastPara.Range.End
$replaceRange = $d.Range($start, $end)
$replaceRange.Text = "Prompt: Design a feature image for `"Buffalo Blitz`" that showcases a happy Maya warrior with glasses in a cartoon style. The image should feature the warrior standing in the midst of the North American prairie, with various animals such as moose, raccoons, pumas, and bears around him. The warrior should be wearing traditional Maya clothing, including a headdress and a necklace made of buffalo bones. The glasses should be modern and stylish, to contrast with the traditional clothing. The image should have a bright and vibrant color scheme, with the warrior smiling while holding a buffalo horn, as if ready to start playing the game. The image should highlight the unique gameplay of Buffalo Blitz with a text overlay that reads `"More Symbols, More Fun: Play Buffalo Blitz Now!`""
